# Manejo de errores e interfaz de usuario
#
# Moves the literal error-message strings (previously typed into column F,
# overlapping the True/False output of column F) into a new column G, and
# puts the actual True/False result back into column F for the "Casos de
# Error" block (rows 20-23). Also appends a new error-case row (24) for an
# invalid "Impuesto de bolsa" value, and documents the expected exception
# text in J21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: "Valor de compra menor a 0" ---
$ws.Range("F20").Value = "False"
$ws.Range("G20").Value = "#ERROR NO EXISTE VALOR DE COMPRA"

# --- Row 21: "Porcentaje de IVA negativo" ---
$ws.Range("F21").Value = "True"
$ws.Range("G21").Value = "#ERROR EL VALOR DEL IVA NO PUEDE SER MENOR A 0"
$ws.Range("J21").Value = "e"

# --- Row 22: "Porcentaje de Impuesto al consumo negativo" ---
$ws.Range("F22").Value = "False"
$ws.Range("G22").Value = "#ERROR EL VALOR DEL IMPUESTO AL CONSUMO NO PUEDE SER MENOR A 0"

# --- Row 23: "Porcentaje de impuesto al licor negativo" ---
$ws.Range("B23").Value = 130000
$ws.Range("F23").Value = "True"
$ws.Range("G23").Value = "#ERROR EL VALOR DEL IMPUESTO AL LICOR NO PUEDE SER MENOR A 0"

# --- Row 24 (new): "Valor invalido de impuesto a bolsa" ---
$ws.Range("A24").Value = "Valor invalido de impuesto a bolsa"
$ws.Range("B24").Value = 23380
$ws.Range("C24").Value = 0.05
$ws.Range("D24").Value = 0.05
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = "None"
$ws.Range("G24").Value = "#ERROR DEBE TENR UN VALOR DE STRING VALIDO"

$ws.Range("D28").Select()
